$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.457.59'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '1.471.49'
$ws.Range('E3').Value = '  +3.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9750'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '275.15'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3647'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.75'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9997'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.456'
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.98'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.161'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001027'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '1.472.31'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9842'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05878'
$ws.Range('E19').Value = '  +3.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.38'
$ws.Range('E20').Value = '  -2.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.452'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.40'
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.94'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').Value = '20.518.32'
$ws.Range('E25').Value = '  +2.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.79'
$ws.Range('E26').Value = '  +6.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.142'
$ws.Range('E27').Value = '  -6.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.25'
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('D29').Value = '1.628.95'
$ws.Range('E29').Value = '  +3.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.57'
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.846'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.975'
$ws.Range('E32').Value = '  -5.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7997'
$ws.Range('E33').Value = '  -2.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07850'
$ws.Range('E34').Value = '  +1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.544'
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05763'
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.159'
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.741'
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.753'
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9750'
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.47'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1872'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5294'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.491'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.98'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.62'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5183'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.771'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06444'
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9879'
$ws.Range('E51').Value = '  -1.26%  '
